# Apply updated crypto price/volume figures (and two name/link row swaps)
# to Sheet1, matching the upstream GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '50.101.53'
$ws.Range("E2").Value = '  +3.46%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '2.675.12'
$ws.Range("E3").Value = '  +6.44%  '

# Row 4: TetherUSD
$ws.Range("E4").Value = '  +0.04%  '

# Row 5: Solana
$ws.Range("D5").Value = '''114.25'
$ws.Range("E5").Value = '  +5.80%  '

# Row 6: BNB
$ws.Range("D6").Value = '''328.85'
$ws.Range("E6").Value = '  +2.50%  '

# Row 7: XRP
$ws.Range("E7").Value = '  +0.65%  '

# Row 8: USDC
$ws.Range("E8").Value = '  +0.03%  '

# Row 9: Cardano
$ws.Range("D9").Value = '''0.561'
$ws.Range("E9").Value = '  +3.40%  '

# Row 10: Avalanche
$ws.Range("D10").Value = '''41.41'
$ws.Range("E10").Value = '  +5.40%  '

# Row 11: Chainlink
$ws.Range("E11").Value = '  +0.66%  '

# Row 12: Dogecoin
$ws.Range("D12").Value = '''0.0827'
$ws.Range("E12").Value = '  +1.95%  '

# Row 13: TRON
$ws.Range("E13").Value = '  +0.78%  '

# Row 14: Polkadot
$ws.Range("E14").Value = '  +3.69%  '

# Row 15: WrappedliquidstakedEther2.0
$ws.Range("D15").Value = '3.085.46'
$ws.Range("E15").Value = '  +6.29%  '

# Row 16: WrappedEther
$ws.Range("D16").Value = '2.681.96'
$ws.Range("E16").Value = '  +6.75%  '

# Row 17: Polygon
$ws.Range("D17").Value = '''0.880'
$ws.Range("E17").Value = '  +5.05%  '

# Row 18: WrappedBTC
$ws.Range("D18").Value = '50.011.62'
$ws.Range("E18").Value = '  +3.63%  '

# Row 19: InternetComputer(DFINITY)
$ws.Range("D19").Value = '''13.29'
$ws.Range("E19").Value = '  +1.26%  '

# Row 20: Uniswap
$ws.Range("D20").Value = '''6.81'
$ws.Range("E20").Value = '  +1.60%  '

# Row 21: ImmutableX
$ws.Range("E21").Value = '  -2.50%  '

# Row 22: ShibaInu
$ws.Range("D22").Value = '0.0₃0968'
$ws.Range("E22").Value = '  +2.57%  '

# Row 23: Litecoin
$ws.Range("D23").Value = '''72.85'
$ws.Range("E23").Value = '  +1.66%  '

# Row 24: BitcoinCash
$ws.Range("D24").Value = '''278.66'
$ws.Range("E24").Value = '  +1.49%  '

# Row 25: PancakeSwap
$ws.Range("E25").Value = '  +2.18%  '

# Row 26: EthereumClassic
$ws.Range("D26").Value = '''26.99'
$ws.Range("E26").Value = '  +3.58%  '

# Row 28: InjectiveProtocol
$ws.Range("D28").Value = '''36.75'
$ws.Range("E28").Value = '  +3.90%  '

# Row 29: Toncoin
$ws.Range("B29").Value = 'Cosmos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D29").Value = '''9.97'
$ws.Range("E29").Value = '  +1.58%  '

# Row 30: Cosmos
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").Value = '''2.22'
$ws.Range("E30").Value = '  +1.13%  '

# Row 31: Kaspa
$ws.Range("D31").Value = '''0.142'
$ws.Range("E31").Value = '  -2.18%  '

# Row 32: OKB
$ws.Range("D32").Value = '''50.53'
$ws.Range("E32").Value = '  +1.71%  '

# Row 33: Celestia
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '''5.46'
$ws.Range("E33").Value = '  +2.19%  '

# Row 34: Filecoin
$ws.Range("B34").Value = 'Celestia'
$ws.Range("C34").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D34").Value = '''19.57'
$ws.Range("E34").Value = '  +0.89%  '

# Row 35: Hedera
$ws.Range("D35").Value = '''0.0802'
$ws.Range("E35").Value = '  +2.00%  '

# Row 36: FirstDigitalUSD
$ws.Range("E36").Value = '  -0.16%  '

# Row 37: ARBITRUM
$ws.Range("D37").Value = '''2.10'
$ws.Range("E37").Value = '  +6.84%  '

# Row 38: RenderToken
$ws.Range("D38").Value = '''4.82'
$ws.Range("E38").Value = '  +3.28%  '

# Row 39: LidoDAOToken
$ws.Range("E39").Value = '  +6.55%  '

# Row 40: Stellar
$ws.Range("E40").Value = '  +1.49%  '

# Row 41: EnergySwap
$ws.Range("B41").Value = 'Monero'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D41").Value = '''124.12'
$ws.Range("E41").Value = '  +3.23%  '

# Row 42: Monero
$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D42").Value = '''22.68'
$ws.Range("E42").Value = '  +2.59%  '

# Row 43: WEMIXToken
$ws.Range("E43").Value = '  +1.61%  '

# Row 44: VeChain
$ws.Range("E44").Value = '  +3.36%  '

# Row 45: NEARProtocol
$ws.Range("E45").Value = '  +3.72%  '

# Row 46: Maker
$ws.Range("D46").Value = '2.083.97'
$ws.Range("E46").Value = '  +3.59%  '

# Row 47: ApeXProtocol
$ws.Range("E47").Value = '  +12.15%  '

# Row 48: Stacks
$ws.Range("D48").Value = '''1.99'
$ws.Range("E48").Value = '  +4.18%  '

# Row 49: FraxShare
$ws.Range("E49").Value = '  +1.76%  '

# Row 50: THORChain
$ws.Range("E50").Value = '  +2.68%  '

# Row 51: BitcoinSV
$ws.Range("D51").Value = '''82.26'
$ws.Range("E51").Value = '  +3.04%  '
